$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("model_attributes")
Write-Host $ws.Name
